$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 697.8461
$ws.Range("I28").Value = 566.65
$ws.Range("K28").Value = 566.65
$ws.Range("M28").Value = -81.64999999999998
$ws.Range("H41").Value = 1395.1428
$ws.Range("J41").Value = 1549.25
$ws.Range("L41").Value = 1549.25
$ws.Range("N41").Value = -2429.25
$ws.Range("H43").Value = 7284.6924
$ws.Range("I43").Value = 6500.5
$ws.Range("J43").Value = 7427.273
$ws.Range("K43").Value = 6500.5
$ws.Range("L43").Value = 7427.273
$ws.Range("M43").Value = -6431.5
$ws.Range("N43").Value = -7565.273
$ws.Range("H51").Value = 7199
$ws.Range("I51").Value = 4499.6665
$ws.Range("J51").Value = 9223.5
$ws.Range("K51").Value = 4499.6665
$ws.Range("L51").Value = 9223.5
$ws.Range("M51").Value = -4015.6665
$ws.Range("N51").Value = -10191.5
$ws.Range("H62").Value = 4039.585
$ws.Range("I62").Value = 3589.4468
$ws.Range("K62").Value = 3589.4468
$ws.Range("M62").Value = -2965.4468
$ws.Range("H65").Value = 4039.585
$ws.Range("I65").Value = 3589.4468
$ws.Range("K65").Value = 17947.234
$ws.Range("M65").Value = -14827.234
$ws.Range("H69").Value = 460999.66
$ws.Range("I69").Value = 998
$ws.Range("J69").Value = 553000
$ws.Range("K69").Value = 2994
$ws.Range("L69").Value = 1659000
$ws.Range("M69").Value = -2120
$ws.Range("N69").Value = -1660748
$ws.Range("H70").Value = 6738.227
$ws.Range("I70").Value = 6216.3335
$ws.Range("J70").Value = 7364.5
$ws.Range("K70").Value = 18649.0005
$ws.Range("L70").Value = 22093.5
$ws.Range("M70").Value = -18379.0005
$ws.Range("N70").Value = -22633.5
$ws.Range("H72").Value = 460999.66
$ws.Range("I72").Value = 998
$ws.Range("J72").Value = 553000
$ws.Range("K72").Value = 8982
$ws.Range("L72").Value = 4977000
$ws.Range("M72").Value = -4614
$ws.Range("N72").Value = -4985736
$ws.Range("H73").Value = 6738.227
$ws.Range("I73").Value = 6216.3335
$ws.Range("J73").Value = 7364.5
$ws.Range("K73").Value = 18649.0005
$ws.Range("L73").Value = 22093.5
$ws.Range("M73").Value = -17713.0005
$ws.Range("N73").Value = -23965.5
$ws.Range("H100").Value = 2782.2
$ws.Range("I100").Value = 1976.25
$ws.Range("J100").Value = 6006
$ws.Range("K100").Value = 1976.25
$ws.Range("L100").Value = 6006
$ws.Range("M100").Value = -1435.25
$ws.Range("N100").Value = -7088
$ws.Range("H103").Value = 1137.1666
$ws.Range("I103").Value = 2066.3333
$ws.Range("J103").Value = 827.44446
$ws.Range("K103").Value = 6198.999899999999
$ws.Range("L103").Value = 2482.33338
$ws.Range("M103").Value = -5612.999899999999
$ws.Range("N103").Value = -3654.33338
$ws.Range("H137").Value = 12986.311
$ws.Range("J137").Value = 2728.5715
$ws.Range("L137").Value = 8185.7145
$ws.Range("N137").Value = -13285.7145

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 5710.8335
$ws.Range("J88").Value = 7721.75
$ws.Range("L88").Value = 7721.75
$ws.Range("N88").Value = -8533.75
$ws.Range("H91").Value = 5710.8335
$ws.Range("J91").Value = 7721.75
$ws.Range("L91").Value = 7721.75
$ws.Range("N91").Value = -10529.75
$ws.Range("H122").Value = 1585.9333
$ws.Range("I122").Value = 1585.9333
$ws.Range("K122").Value = 4757.7999
$ws.Range("M122").Value = -2307.7999
$ws.Range("H124").Value = 25365.6
$ws.Range("J124").Value = 25365.6
$ws.Range("L124").Value = 25365.6
$ws.Range("N124").Value = -35185.6
$ws.Range("H132").Value = 1406.375
$ws.Range("I132").Value = 1145.7646
$ws.Range("J132").Value = 2883.1667
$ws.Range("K132").Value = 3437.2938
$ws.Range("L132").Value = 8649.500100000001
$ws.Range("M132").Value = -907.2937999999999
$ws.Range("N132").Value = -13709.5001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 12833.259
$ws.Range("I20").Value = 18666
$ws.Range("J20").Value = 1167.7778
$ws.Range("K20").Value = 18666
$ws.Range("L20").Value = 1167.7778
$ws.Range("M20").Value = -18419
$ws.Range("N20").Value = -1661.7778
$ws.Range("H86").Value = 1452.9429
$ws.Range("I86").Value = 1403.5186
$ws.Range("K86").Value = 1403.5186
$ws.Range("M86").Value = -280.5186000000001
$ws.Range("H89").Value = 1452.9429
$ws.Range("I89").Value = 1403.5186
$ws.Range("K89").Value = 7017.593000000001
$ws.Range("M89").Value = -1401.593000000001
$ws.Range("H105").Value = 3881.4167
$ws.Range("I105").Value = 3018.3333
$ws.Range("J105").Value = 4744.5
$ws.Range("K105").Value = 3018.3333
$ws.Range("L105").Value = 4744.5
$ws.Range("M105").Value = -1271.3333
$ws.Range("N105").Value = -8238.5
$ws.Range("H134").Value = 1734.3489
$ws.Range("I134").Value = 1308.9714
$ws.Range("K134").Value = 3926.9142
$ws.Range("M134").Value = -1391.9142

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 8351
$ws.Range("I2").Value = 8351
$ws.Range("K2").Value = 8351
$ws.Range("M2").Value = -8238
$ws.Range("H5").Value = 633.4286
$ws.Range("I5").Value = 506
$ws.Range("J5").Value = 803.3333
$ws.Range("K5").Value = 506
$ws.Range("L5").Value = 803.3333
$ws.Range("M5").Value = -394
$ws.Range("N5").Value = -1027.3333
$ws.Range("H8").Value = 3109
$ws.Range("I8").Value = 5998
$ws.Range("K8").Value = 5998
$ws.Range("M8").Value = -5858
$ws.Range("H10").Value = 426.25
$ws.Range("I10").Value = 426.25
$ws.Range("K10").Value = 426.25
$ws.Range("M10").Value = -287.25
$ws.Range("H12").Value = 141632290
$ws.Range("I12").Value = 225156260
$ws.Range("J12").Value = 30267000
$ws.Range("K12").Value = 225156260
$ws.Range("L12").Value = 30267000
$ws.Range("M12").Value = -225156090
$ws.Range("N12").Value = -30267340
$ws.Range("H14").Value = 17048.5
$ws.Range("J14").Value = 21399.666
$ws.Range("L14").Value = 21399.666
$ws.Range("N14").Value = -21739.666
$ws.Range("H15").Value = 586.4
$ws.Range("I15").Value = 625
$ws.Range("K15").Value = 625
$ws.Range("M15").Value = -455
$ws.Range("H86").Value = 15555.632
$ws.Range("I86").Value = 19239.3
$ws.Range("J86").Value = 11462.667
$ws.Range("K86").Value = 19239.3
$ws.Range("L86").Value = 11462.667
$ws.Range("M86").Value = -18116.3
$ws.Range("N86").Value = -13708.667
$ws.Range("H89").Value = 15555.632
$ws.Range("I89").Value = 19239.3
$ws.Range("J89").Value = 11462.667
$ws.Range("K89").Value = 96196.5
$ws.Range("L89").Value = 57313.335
$ws.Range("M89").Value = -90580.5
$ws.Range("N89").Value = -68545.33499999999
$ws.Range("H100").Value = 84365.57000000001
$ws.Range("J100").Value = 84365.57000000001
$ws.Range("L100").Value = 84365.57000000001
$ws.Range("N100").Value = -86529.57000000001
$ws.Range("H122").Value = 1508.5312
$ws.Range("I122").Value = 1460.5
$ws.Range("J122").Value = 1652.625
$ws.Range("K122").Value = 4381.5
$ws.Range("L122").Value = 4957.875
$ws.Range("M122").Value = -1931.5
$ws.Range("N122").Value = -9857.875

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 8766.556
$ws.Range("I116").Value = 2599.3333
$ws.Range("K116").Value = 7797.999899999999
$ws.Range("M116").Value = -4355.999899999999
$ws.Range("H118").Value = 3000
$ws.Range("J118").Value = 3000
$ws.Range("L118").Value = 9000
$ws.Range("N118").Value = -11486
$ws.Range("H136").Value = 2694
$ws.Range("I136").Value = 2392.8
$ws.Range("K136").Value = 7178.400000000001
$ws.Range("M136").Value = -2078.400000000001
$ws.Range("H139").Value = 1780
$ws.Range("I139").Value = 1616.1765
$ws.Range("J139").Value = 2708.3333
$ws.Range("K139").Value = 4848.529500000001
$ws.Range("L139").Value = 8124.999899999999
$ws.Range("M139").Value = 291.4704999999994
$ws.Range("N139").Value = -18404.9999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 192.94737
$ws.Range("I2").Value = 167
$ws.Range("J2").Value = 211.81818
$ws.Range("K2").Value = 167
$ws.Range("L2").Value = 211.81818
$ws.Range("M2").Value = -54
$ws.Range("N2").Value = -437.81818
$ws.Range("H70").Value = 7493.067
$ws.Range("I70").Value = 6950.1055
$ws.Range("J70").Value = 8430.909
$ws.Range("K70").Value = 6950.1055
$ws.Range("L70").Value = 8430.909
$ws.Range("M70").Value = -6680.1055
$ws.Range("N70").Value = -8970.909
$ws.Range("H73").Value = 7493.067
$ws.Range("I73").Value = 6950.1055
$ws.Range("J73").Value = 8430.909
$ws.Range("K73").Value = 6950.1055
$ws.Range("L73").Value = 8430.909
$ws.Range("M73").Value = -6014.1055
$ws.Range("N73").Value = -10302.909
$ws.Range("H97").Value = 1785.4814
$ws.Range("I97").Value = 1219.625
$ws.Range("K97").Value = 1219.625
$ws.Range("M97").Value = -723.625

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1400.75
$ws.Range("I22").Value = 1247.3846
$ws.Range("J22").Value = 1582
$ws.Range("K22").Value = 1247.3846
$ws.Range("L22").Value = 1582
$ws.Range("M22").Value = -952.3846000000001
$ws.Range("N22").Value = -2172
$ws.Range("H27").Value = 1400.75
$ws.Range("I27").Value = 1247.3846
$ws.Range("J27").Value = 1582
$ws.Range("K27").Value = 1247.3846
$ws.Range("L27").Value = 1582
$ws.Range("M27").Value = -1140.3846
$ws.Range("N27").Value = -1796
$ws.Range("H40").Value = 2108.3447
$ws.Range("I40").Value = 1989.68
$ws.Range("K40").Value = 1989.68
$ws.Range("M40").Value = -1853.68
$ws.Range("H75").Value = 90000
$ws.Range("J75").Value = 90000
$ws.Range("L75").Value = 90000
$ws.Range("N75").Value = -91872
$ws.Range("H78").Value = 90000
$ws.Range("J78").Value = 90000
$ws.Range("L78").Value = 270000
$ws.Range("N78").Value = -279360

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2301.7778
$ws.Range("I126").Value = 1964.5
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 5893.5
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -3423.5
$ws.Range("N126").Value = -19940
